$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "294.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.07%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.20%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.932"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.93%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07367"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.68%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.295"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "29.27%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.695"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.56%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.758"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.52%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9138"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.06%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1688"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.75%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08371"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "11.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08241"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.88%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.10%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1008"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.77%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001509"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.23%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005776"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.17%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.42%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.35%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.39%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.971"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-7.06%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04547"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.06%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001209"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.35%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004336"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.68%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.81%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.95%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04445"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.02%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007309"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.16%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.008802"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1327"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.58%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002061"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.72%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009107"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.19%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006026"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.53%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.24%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.92%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.24%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.24%"

Write-Host "Updated crypto symbol prices and volumes."